$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.842.46"
$ws.Range("E2").Value = "  +1.07%  "

$ws.Range("D3").Value = "'3.123.64"
$ws.Range("E3").Value = "  +0.75%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'532.53"
$ws.Range("E5").Value = "  +1.47%  "

$ws.Range("D6").Value = "'138.56"
$ws.Range("E6").Value = "  +0.80%  "

$ws.Range("E8").Value = "  +6.40%  "

$ws.Range("E9").Value = "  +0.15%  "

$ws.Range("E10").Value = "  +0.43%  "

$ws.Range("D11").Value = "'0.413"
$ws.Range("E11").Value = "  +2.96%  "

$ws.Range("D12").Value = "'0.139"
$ws.Range("E12").Value = "  +3.13%  "

$ws.Range("D13").Value = "'3.665.16"
$ws.Range("E13").Value = "  +0.92%  "

$ws.Range("D14").Value = "'25.77"
$ws.Range("E14").Value = "  +0.96%  "

$ws.Range("E15").Value = "  +1.81%  "

$ws.Range("D16").Value = "'57.933.98"
$ws.Range("E16").Value = "  +1.04%  "

$ws.Range("D17").Value = "'3.136.98"
$ws.Range("E17").Value = "  +1.19%  "

$ws.Range("E18").Value = "  +2.79%  "

$ws.Range("D19").Value = "'12.78"
$ws.Range("E19").Value = "  +2.03%  "

$ws.Range("D20").Value = "'8.12"
$ws.Range("E20").Value = "  +2.57%  "

$ws.Range("D21").Value = "'373.71"
$ws.Range("E21").Value = "  +7.19%  "

$ws.Range("D22").Value = "'0.998"
$ws.Range("E22").Value = "  -0.10%  "

$ws.Range("E23").Value = "  -1.34%  "

$ws.Range("D24").Value = "'69.52"
$ws.Range("E24").Value = "  +2.17%  "

$ws.Range("E25").Value = "  +1.02%  "

$ws.Range("E26").Value = "  +0.21%  "

$ws.Range("E27").Value = "  +0.30%  "

$ws.Range("D28").Value = "'0.0₃0874"
$ws.Range("E28").Value = "  -1.75%  "

$ws.Range("D29").Value = "'7.41"
$ws.Range("E29").Value = "  +0.88%  "

$ws.Range("D30").Value = "'6.16"
$ws.Range("E30").Value = "  +2.90%  "

$ws.Range("E31").Value = "  -0.28%  "

$ws.Range("D32").Value = "'21.46"
$ws.Range("E32").Value = "  +2.95%  "

$ws.Range("D33").Value = "'5.13"
$ws.Range("E33").Value = "  +2.25%  "

$ws.Range("E34").Value = "  +2.80%  "

$ws.Range("D35").Value = "'160.21"
$ws.Range("E35").Value = "  +0.72%  "

$ws.Range("D36").Value = "'6.15"
$ws.Range("E36").Value = "  +1.47%  "

$ws.Range("E37").Value = "  +3.52%  "

$ws.Range("E38").Value = "  -2.83%  "

$ws.Range("E39").Value = "  +4.58%  "

$ws.Range("D40").Value = "'0.0671"
$ws.Range("E40").Value = "  +2.10%  "

$ws.Range("D41").Value = "'2.535.26"
$ws.Range("E41").Value = "  +5.45%  "

$ws.Range("D42").Value = "'4.09"
$ws.Range("E42").Value = "  +0.14%  "

$ws.Range("E43").Value = "  +0.71%  "

$ws.Range("D44").Value = "'37.86"
$ws.Range("E44").Value = "  +3.29%  "

$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0270"
$ws.Range("E45").Value = "  +2.21%  "

$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").Value = "'1.00"
$ws.Range("E46").Value = "  -0.06%  "

$ws.Range("D47").Value = "'0.978"
$ws.Range("E47").Value = "  +0.62%  "

$ws.Range("D48").Value = "'6.15"
$ws.Range("E48").Value = "  +2.63%  "

$ws.Range("D49").Value = "'19.73"
$ws.Range("E49").Value = "  +0.06%  "

$ws.Range("D50").Value = "'0.749"
$ws.Range("E50").Value = "  -1.85%  "

$ws.Range("E51").Value = "  +3.17%  "
